# Fruta / hortaliza, semanal
# Insert two new weekly data rows (rows 37 and 38) into the Espárragos sheet.
# This shifts the existing rows 37-115 down to 39-117 and populates the
# two newly-inserted rows with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 37 (pushes old row 37 -> 39, ..., old row 115 -> 117)
$ws.Rows.Item(37).EntireRow.Insert()
$ws.Rows.Item(37).EntireRow.Insert()

# ---- Row 37 (new data) ----
$ws.Range("A37").Value = 6
$ws.Range("B37").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C37").Value = "Metropolitana"
$ws.Range("D37").Value = 44525
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = 300000000
$ws.Range("G37").Value = "Espárragos"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Banquete"
$ws.Range("J37").Value = 350
$ws.Range("K37").Value = 1400
$ws.Range("L37").Value = 1500
$ws.Range("M37").Value = 1466
$ws.Range("N37").Value = "$/kilo"
$ws.Range("O37").Value = "Provincia de Linares"
$ws.Range("P37").Value = 1466
$ws.Range("Q37").Value = 1
$ws.Range("R37").Value = "Hortaliza"

# ---- Row 38 (new data) ----
$ws.Range("A38").Value = 6
$ws.Range("B38").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C38").Value = "Metropolitana"
$ws.Range("D38").Value = 44525
$ws.Range("E38").Value = 13
$ws.Range("F38").Value = 300000000
$ws.Range("G38").Value = "Espárragos"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 880
$ws.Range("K38").Value = 1000
$ws.Range("L38").Value = 1300
$ws.Range("M38").Value = 1205
$ws.Range("N38").Value = "$/kilo"
$ws.Range("O38").Value = "Provincia de Linares"
$ws.Range("P38").Value = 1205
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = "Hortaliza"
